$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Saldo (balance) value for account 001882235 / LAGO (row 2)
$ws.Range("C2").Value = 229707.72

# Delete the entire data row for account 004332544 / CELIA (row 14);
# all subsequent rows shift up by one.
$ws.Rows(14).Delete()
